$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 1) "Venom" paragraph: drop the spellcheck proofErr markers and merge the
#    "Venom" run + the trailing space run into a single run "Venom "
$d.Paragraphs(2).Range.InsertXML("<w:p $wNs><w:pPr><w:pStyle w:val=`"PargrafodaLista`"/></w:pPr><w:r><w:t xml:space=`"preserve`">Venom </w:t></w:r></w:p>") | Out-Null

# 2) "Sherek" paragraph: drop the spellcheck proofErr markers around it
$d.Paragraphs(3).Range.InsertXML("<w:p $wNs><w:pPr><w:pStyle w:val=`"PargrafodaLista`"/></w:pPr><w:r><w:t>Sherek</w:t></w:r></w:p>") | Out-Null

# 3) "Moana" paragraph: drop the spellcheck proofErr markers around it
$d.Paragraphs(4).Range.InsertXML("<w:p $wNs><w:pPr><w:pStyle w:val=`"PargrafodaLista`"/></w:pPr><w:r><w:t>Moana</w:t></w:r></w:p>") | Out-Null

# 4) Add the new movie "Jackss" into the previously empty 7th paragraph
$d.Paragraphs(7).Range.Text = "Jackss"
